$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (aggiornamento fino a 9 agosto 2021)
# Columns: row, A (date serial), B (nuovi pos.), C (somma mobile 7gg.), D (somma mobile 7gg. per 100mila abitanti)
$data = @(
  @(329, 44403, 1, 8, 129.3870289503477),
  @(330, 44404, 0, 8, 129.3870289503477),
  @(331, 44405, 0, 8, 129.3870289503477),
  @(332, 44406, 0, 8, 129.3870289503477),
  @(333, 44407, 0, 6, 97.0402717127608),
  @(334, 44408, 0, 3, 48.5201358563804),
  @(335, 44409, 1, 2, 32.34675723758694),
  @(336, 44410, 0, 1, 16.17337861879347),
  @(337, 44411, 0, 1, 16.17337861879347),
  @(338, 44412, 0, 1, 16.17337861879347),
  @(339, 44413, 2, 3, 48.5201358563804),
  @(340, 44414, 0, 3, 48.5201358563804),
  @(341, 44415, 1, 4, 64.69351447517387),
  @(342, 44416, 0, 3, 48.5201358563804),
  @(343, 44417, 0, 3, 48.5201358563804)
)

$lastExistingRow = 328

foreach ($row in $data) {
    $r = $row[0]

    # Copy formatting (style) from the last row of the existing table so the
    # new row matches the look of previous rows (date style with borders etc.)
    $srcA = $ws.Cells.Item($lastExistingRow, 1)
    $srcB = $ws.Cells.Item($lastExistingRow, 2)
    $srcC = $ws.Cells.Item($lastExistingRow, 3)
    $srcD = $ws.Cells.Item($lastExistingRow, 4)

    $dstA = $ws.Cells.Item($r, 1)
    $dstB = $ws.Cells.Item($r, 2)
    $dstC = $ws.Cells.Item($r, 3)
    $dstD = $ws.Cells.Item($r, 4)

    $srcA.Copy($dstA)
    $srcB.Copy($dstB)
    $srcC.Copy($dstC)
    $srcD.Copy($dstD)

    $dstA.Value2 = $row[1]
    $dstB.Value2 = $row[2]
    $dstC.Value2 = $row[3]
    $dstD.Value2 = $row[4]
}
